$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "35.431 * (3.874)"
$ws.Range("C2").Value = "0.02 * (0.0034)"

$ws.Range("B3").Value = "59.625 * (7.169)"
$ws.Range("C3").Value = "0.0395 * (0.006)"

$ws.Range("B4").Value = "77.056 * (12.697)"
$ws.Range("C4").Value = "0.0402 * (0.0142)"

$ws.Range("B5").Value = "28.224 * (9.213)"
$ws.Range("C5").Value = "0.0237 * (0.0069)"

$ws.Range("B6").Value = "24.683  (9.993)"
$ws.Range("C6").Value = "-0.0034  (0.0192)"

$ws.Range("B7").Value = "5.742  (7.11)"
$ws.Range("C7").Value = "0.0007  (0.0073)"

$ws.Range("B8").Value = "15.625  (8.773)"
$ws.Range("C8").Value = "0.0042  (0.0047)"
